# feat: add 2022-Q1 data
#
# 1. Insert a new sheet "2022-Q1" (fund-level detail) right before the
#    "总计" (total) summary sheet, duplicating the layout/format of an
#    existing quarter sheet ("2021-Q4") so header/column-A styling matches.
# 2. Populate that new sheet with the 2022-Q1 fund holdings.
# 3. Prepend a "2022-Q1" row to the "总计" summary sheet, pushing the
#    existing rows down and renumbering the index column.

$wb = $excel.ActiveWorkbook

$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")

# --- 1. Create the new "2022-Q1" sheet by copying "2021-Q4", which already
#        carries the right sheet/header/column styling, then drop it right
#        before "总计" and rename it.
#
# NOTE: worksheet variables appear to be resolved by position, not a live
# object binding, so `$totalSheet` must be re-fetched by name after the
# sheet collection changes below (its index shifts from 4 to 5).
$q4Sheet.Copy($totalSheet, $null)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"
$totalSheet = $wb.Worksheets.Item("总计")

# --- 2. Overwrite the header row (row 1) text.
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# Make sure the fund code (column B, has meaningful leading zeros) and the
# numeric-looking fund figures (columns D-G) are stored as text, matching
# the source data convention used by every other quarter sheet.
$newSheet.Range("B2:B7").NumberFormat = "@"
$newSheet.Range("D2:G7").NumberFormat = "@"

# Row 2
$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "163302"
$newSheet.Cells.Item(2,3).Value = "大摩资源优选混合(LOF)"
$newSheet.Cells.Item(2,4).Value = "5.82"
$newSheet.Cells.Item(2,5).Value = "81.78"
$newSheet.Cells.Item(2,6).Value = "3.58"
$newSheet.Cells.Item(2,7).Value = "0.2084"
$newSheet.Cells.Item(2,8).Value = 7

# Row 3
$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).Value = "008347"
$newSheet.Cells.Item(3,3).Value = "中信建投价值甄选混合A"
$newSheet.Cells.Item(3,4).Value = "4.45"
$newSheet.Cells.Item(3,5).Value = "72.24"
$newSheet.Cells.Item(3,6).Value = "2.02"
$newSheet.Cells.Item(3,7).Value = "0.0899"
$newSheet.Cells.Item(3,8).Value = 10

# Row 4
$newSheet.Cells.Item(4,1).Value = 2
$newSheet.Cells.Item(4,2).Value = "003822"
$newSheet.Cells.Item(4,3).Value = "中信建投行业轮换混合A"
$newSheet.Cells.Item(4,4).Value = "3.07"
$newSheet.Cells.Item(4,5).Value = "72.09"
$newSheet.Cells.Item(4,6).Value = "2.24"
$newSheet.Cells.Item(4,7).Value = "0.0688"
$newSheet.Cells.Item(4,8).Value = 10

# Row 5
$newSheet.Cells.Item(5,1).Value = 3
$newSheet.Cells.Item(5,2).Value = "007468"
$newSheet.Cells.Item(5,3).Value = "中信建投策略精选混合A"
$newSheet.Cells.Item(5,4).Value = "0.94"
$newSheet.Cells.Item(5,5).Value = "78.13"
$newSheet.Cells.Item(5,6).Value = "2.80"
$newSheet.Cells.Item(5,7).Value = "0.0263"
$newSheet.Cells.Item(5,8).Value = 8

# Row 6
$newSheet.Cells.Item(6,1).Value = 4
$newSheet.Cells.Item(6,2).Value = "003823"
$newSheet.Cells.Item(6,3).Value = "中信建投行业轮换混合C"
$newSheet.Cells.Item(6,4).Value = "0.64"
$newSheet.Cells.Item(6,5).Value = "72.09"
$newSheet.Cells.Item(6,6).Value = "2.24"
$newSheet.Cells.Item(6,7).Value = "0.0143"
$newSheet.Cells.Item(6,8).Value = 10

# Row 7
$newSheet.Cells.Item(7,1).Value = 5
$newSheet.Cells.Item(7,2).Value = "007469"
$newSheet.Cells.Item(7,3).Value = "中信建投策略精选混合C"
$newSheet.Cells.Item(7,4).Value = "0.40"
$newSheet.Cells.Item(7,5).Value = "78.13"
$newSheet.Cells.Item(7,6).Value = "2.80"
$newSheet.Cells.Item(7,7).Value = "0.0112"
$newSheet.Cells.Item(7,8).Value = 8

# --- 3. Prepend the 2022-Q1 summary row to "总计", shifting the existing
#        rows down by one and renumbering the index column (A).

# Extend row formatting down to the new last row (5) by copying row 4's
# format before overwriting values, so every row keeps the same column-A
# style used throughout the sheet.
$totalSheet.Range("A4:D4").Copy()
$totalSheet.Range("A5:D5").PasteSpecial(-4122)

# Row 5 <- old row 4 (2021-Q1)
$totalSheet.Cells.Item(5,1).Value = 3
$totalSheet.Cells.Item(5,2).Value = "2021-Q1"
$totalSheet.Cells.Item(5,3).Value = 2
$totalSheet.Cells.Item(5,4).Value = 0.08

# Row 4 <- old row 3 (2021-Q2)
$totalSheet.Cells.Item(4,1).Value = 2
$totalSheet.Cells.Item(4,2).Value = "2021-Q2"
$totalSheet.Cells.Item(4,3).Value = 2
$totalSheet.Cells.Item(4,4).Value = 0.08

# Row 3 <- old row 2 (2021-Q4)
$totalSheet.Cells.Item(3,1).Value = 1
$totalSheet.Cells.Item(3,2).Value = "2021-Q4"
$totalSheet.Cells.Item(3,3).Value = 6
$totalSheet.Cells.Item(3,4).Value = 0.43

# Row 2 <- new 2022-Q1 row
$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q1"
$totalSheet.Cells.Item(2,3).Value = 6
$totalSheet.Cells.Item(2,4).Value = 0.42
